$wb = $excel.ActiveWorkbook

# Work on the "PE, SEPTEMBER" worksheet (second sheet, sheet2.xml / rId2)
$ws = $wb.Worksheets.Item("PE, SEPTEMBER")

# Update selection (I9 -> I10) to match final state
$ws.Activate()
$ws.Range("I10").Select()

# Row 9: add date, OR/invoice number, and gross amount formula
$ws.Range("C9").Value = 45937
$ws.Range("G9").Value = 517976202
$ws.Range("I9").Formula = "=1353132-56255.04"

# Row 40: clear out the CARL'S VULCANIZING SHOP entry entirely
$ws.Range("C40:I40").ClearContents()

$wb.Save()
